$d = $word.ActiveDocument

# --- Step 0: header/footer distance becomes 0 in the saved section properties ---
$ps = $d.Sections.Item(1).PageSetup
$ps.HeaderDistance = 0
$ps.FooterDistance = 0

# --- Step 1: insert two new paragraphs before the "Cuando salen tres Gy:" paragraph ---
# Before: P1(empty) P2(Cuando...) P3(pic) P4(empty) P5(Si...)
$pCuando = $d.Paragraphs.Item(2)
$pCuando.Range.InsertParagraphBefore()
$pCuando.Range.InsertParagraphBefore()
# Now: P1(empty) P2(empty-NEW) P3(empty-NEW) P4(Cuando...) P5(pic) P6(empty) P7(Si...)

$p3 = $d.Paragraphs.Item(3)
$p3.Range.Text = "1. las fechas deberian comenzar en las variables 114 y 120?"
# Now P3 has the new question 1 text

# --- Step 2: replace the "Cuando salen tres Gy:" paragraph (index 4) cleanly ---
# (delete + recreate to drop leftover proofErr spell-check markers)
$pCuando = $d.Paragraphs.Item(4)
$pCuando.Range.InsertParagraphBefore()
$newQ2 = $d.Paragraphs.Item(4)
$newQ2.Range.Text = "2. Cuando salen tres Gy, "
$d.Paragraphs.Item(5).Range.Delete()
# Now: P1(empty) P2(empty) P3("1. ...") P4("2. Cuando salen tres Gy, ") P5(pic) P6(empty) P7(Si...)

# --- Step 3: replace the "Si el paciente..." paragraph (index 7) with expanded text ---
$pSi = $d.Paragraphs.Item(7)
$pSi.Range.InsertParagraphBefore()
$newQ3 = $d.Paragraphs.Item(7)
$newQ3.Range.Text = "3. Si el paciente está recibiendo radioterapia, pero aún no termina, 57402094, que se coloca en la fecha de finalización y solo se llena 114?"
$d.Paragraphs.Item(8).Range.Delete()
# Now: P1(empty) P2(empty) P3("1. ...") P4("2. ...") P5(pic) P6(empty) P7("3. ...")

# --- Step 4: append two new empty paragraphs at the end of the body text (after P7) ---
$pSi = $d.Paragraphs.Item(7)
$pSi.Range.InsertParagraphAfter()
$pSi.Range.InsertParagraphAfter()
# Now: P1..P7 as above, P8(empty-NEW), P9(empty-NEW)

# last paragraph carries explicit spacing (before=0, after=8pt -> w:after=160)
$pLast = $d.Paragraphs.Item($d.Paragraphs.Count)
$pLast.Range.ParagraphFormat.SpaceBefore = 0
$pLast.Range.ParagraphFormat.SpaceAfter = 8
